# Add a new "ticklabels" function row to the summary table on the
# "table_of_functions.csv" sheet, right after the "ax" row (under the
# "visualization" category) and before the "2D density" category block.
#
# This mirrors inserting a whole row at row 17 (pushing every row from the
# old row 17 onward down by one) and filling in Function/Description for
# the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table_of_functions.csv")

# Insert a new blank row at row 17 - everything below shifts down one row.
$ws.Rows.Item(17).Insert()

# Populate the new row: no category (continues the "visualization" group
# started above), function name in column B, description in column C.
$ws.Cells.Item(17, 2).Value = "ticklabels"
$ws.Cells.Item(17, 3).Value = "Label powers of 10 tickmarks"

# The freshly inserted row inherits formatting from the row above it for
# every column, but column C in the rest of the table uses a different
# style (s=7) than column B (s=6). Copy that formatting down from the row
# below (the old row 17, now row 18) onto the new C17 so it matches the
# rest of the table.
$ws.Cells.Item(18, 3).Copy()
$ws.Cells.Item(17, 3).PasteSpecial(-4122)

# Match the author's final selection recorded in the sheet.
[void]$ws.Cells.Item(18, 2).Select()
